$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fota_page")

# Row 5: rename the firmware state message id element (drop the "_idle" suffix
# so it now refers to the generic firmware_state_message element, used to
# support multiple devices).
$ws.Range("B5").Value = "com.tcl.fota.system:id/firmware_state_message"

# Rows 6-9 (new_version/2BBF, old_version/2BBE, launcher_package/
# com.tcl.android.launcher, expect_update_time/5) are unchanged.

# Add two new elements needed for the multi-device update flow.
$ws.Range("A10").Value = "permission_agree"
$ws.Range("B10").Value = "AGREE"
$ws.Range("C10").Value = "str"

$ws.Range("A11").Value = "update_restart_time"
$ws.Range("B11").Value = 600
$ws.Range("C11").Value = "str"

# Move the active selection to B7, matching the final sheet state.
$ws.Range("B7").Select()
